# Add payment 09876543 (Cash) 2025-08-18T18:03:44
#
# The existing row 53 (phone "09876543" stored as text with a leading
# zero) is corrected to a plain number (09876543 -> 9876543, no leading
# zero) and a brand-new row 54 is appended carrying the original
# leading-zero phone text together with the rest of the payment record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 53, column A: text "09876543" -> number 9876543 ---------------
$ws.Cells.Item(53, 1).Value = 9876543

# --- Append new row 54: the actual new payment ------------------------------
# Column A keeps the leading zero, so it must stay text (leading apostrophe
# forces text entry, exactly like typing '09876543 into the cell in Excel).
$ws.Cells.Item(54, 1).Value = "'09876543"

$ws.Cells.Item(54, 2).Value = ""
$ws.Cells.Item(54, 3).Value = "Cash"
$ws.Cells.Item(54, 4).Value = "2025-08-18T18:03:44"
$ws.Cells.Item(54, 5).Value = 120
$ws.Cells.Item(54, 6).Value = ""
$ws.Cells.Item(54, 7).Value = 120
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 0
